$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Rows 2-13: Runmode column C flips from "N" to "Y"
$ws.Range("C2:C13").Value = "Y"

# Row 14: Runmode column C flips from "Y" to "N"
$ws.Range("C14").Value = "N"

# Rows 15-16: Runmode column C flips from "N" to "Y"; also normalize their
# cell style (border-only, no explicit fill) to match the rest of the column.
$ws.Range("C2").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C15").Value = "Y"
$ws.Range("C16").Value = "Y"

# Move the active selection from C17 to C14
[void]$ws.Range("C14").Select()
